# Update post last call, with all examples and images
#
# 1) The "datetimeFigureOut" date field (08/06/2020 -> 12/06/2020) that is
#    reproduced on the slide master and on every slide layout's date
#    placeholder.
# 2) The "specimenRequirements" label on slide 1 becomes "specimenRequested".

$p = $ppt.ActivePresentation

# --- 1. Fix the date placeholder text everywhere it appears -----------------
$OldDate = "08/06/2020"
$NewDate = "12/06/2020"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $OldDate) {
                $tr.Text = $NewDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout hanging off the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2. Rename "specimenRequirements" to "specimenRequested" on slide 1 -----
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shape = $shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)
            if ($para.Text.TrimEnd([char]13) -eq "   specimenRequirements ") {
                $para.Text = "   specimenRequested "
            }
        }
    }
}
